# Mise à jour des références sur le relevé d'erreur
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Référence" (F) column links ---------------------------
# Order matters: it controls the order new strings are appended to the
# shared-string table (F12 first, then F2, F5, F4, F6, F7, F8, F9, F10,
# F11, F13, F18 — matching the authoring order in the target workbook).

$ws.Range("F12").Value = "https://www.orixa-media.com/academie/referencement-naturel/technique/balise-meta-description/?utm_source=google&utm_medium=orixa-site&gclid=CjwKCAjwhaaKBhBcEiwA8acsHGvULoJNUxZxReWYFnntQHjFM_BnKldrTaUi-9S0kcUfg1C5bt3uSxoCO1EQAvD_BwE"
$ws.Range("F2").Value  = "https://smartkeyword.io/seo-on-page-balise-title/"
$ws.Range("F5").Value  = "https://black.bird.eu/fr/blog/performances-accelerer-le-chargement-des-images-sur-magento-2.html"
$ws.Range("F4").Value  = "https://blog.hubspot.fr/marketing/reduire-duree-chargement-page-web"
$ws.Range("F6").Value  = "https://fr.semrush.com/blog/texte-alternatif/"
$ws.Range("F7").Value  = "https://www.tech-wiki.online/fr/javascript-async-defer.html"
$ws.Range("F8").Value  = "https://maxime-benard.fr/article/42/fichier-htaccess-a-quoi-ca-sert"
$ws.Range("F9").Value  = "https://minifier.org/"
$ws.Range("F10").Value = "https://fr.ryte.com/magazine/utiliser-mise-cache-navigateur-accelerer-site-web"
$ws.Range("F11").Value = "http://www.pompage.net/traduction/Bien-utiliser-le-texte-alternatif"
$ws.Range("F13").Value = "https://developer.mozilla.org/fr/docs/Web/HTML/Global_attributes/lang"
$ws.Range("F18").Value = "https://wbcreation.fr/normes-w3c.html"

# Cells that reuse strings already introduced above
$ws.Range("F14").Value = "https://www.orixa-media.com/academie/referencement-naturel/technique/balise-meta-description/?utm_source=google&utm_medium=orixa-site&gclid=CjwKCAjwhaaKBhBcEiwA8acsHGvULoJNUxZxReWYFnntQHjFM_BnKldrTaUi-9S0kcUfg1C5bt3uSxoCO1EQAvD_BwE"
$ws.Range("F15").Value = "https://www.orixa-media.com/academie/referencement-naturel/technique/balise-meta-description/?utm_source=google&utm_medium=orixa-site&gclid=CjwKCAjwhaaKBhBcEiwA8acsHGvULoJNUxZxReWYFnntQHjFM_BnKldrTaUi-9S0kcUfg1C5bt3uSxoCO1EQAvD_BwE"
$ws.Range("F16").Value = "http://www.pompage.net/traduction/Bien-utiliser-le-texte-alternatif"

# F3 and F17 keep their original reference (MDN mobile checklist) untouched.

# --- Rebuild the hyperlinks collection -----------------------------------
# Drop the old hyperlinks (F2, F3:F5, F6:F18) and keep just the two that
# remain in the updated sheet: F2 (unchanged target) and a new one on F12.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://developer.mozilla.org/fr/docs/Accessibilit%C3%A9/Checklist_accessibilite_mobile") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.orixa-media.com/academie/referencement-naturel/technique/balise-meta-description/?utm_source=google&utm_medium=orixa-site&gclid=CjwKCAjwhaaKBhBcEiwA8acsHGvULoJNUxZxReWYFnntQHjFM_BnKldrTaUi-9S0kcUfg1C5bt3uSxoCO1EQAvD_BwE") | Out-Null

# Hyperlinks.Add() re-stamps the cell style (adds an explicit applyFont
# flag), so restore the original "Lien hypertexte" cell style that was
# already in place for these two cells.
$ws.Range("F2").Style = "Lien hypertexte"
$ws.Range("F12").Style = "Lien hypertexte"

# --- Row height tweak on the first data row -------------------------------
$ws.Rows(2).RowHeight = 18.75

# --- Move the active cell / selection -------------------------------------
$ws.Range("F25").Select()
